$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold numeric-looking text; force text format
# so Excel COM does not reinterpret them as numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '256.60'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.18%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.07'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-3.30%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.540'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-13.01%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05902'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.40%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.626'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.46%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8600'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-1.92%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9349'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-5.33%'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.62%'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.03629'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.99%'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07082'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.97%'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03231'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.98%'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09203'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.40%'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001546'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.09%'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006051'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.66%'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006110'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.30%'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.60%'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.196'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.99%'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.16%'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3056'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.12%'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.03%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.851'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '9.06%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04221'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.57%'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.50%'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-6.00%'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.22%'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.06%'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.44%'

$ws.Range("B41").Value = 'BKEXToken'

$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1101'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.45%'

$ws.Range("B42").Value = 'KickToken'

$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003950'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-28.02%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002410'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '1.57%'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '6.10%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005458'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.74%'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.06021'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-29.51%'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06870'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '3,107.01%'
